# Update the Start Date (row 4) and End Date (row 5) values for columns B, C, D
# on the "100_1" worksheet. The cells already carry a date number format, so
# simply assigning the new serial date values preserves formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("100_1")

# Row 4 - Start Date: 32511 (1/3/1989) -> 31783 (1/6/1987)
$ws.Range("B4").Value = 31783
$ws.Range("C4").Value = 31783
$ws.Range("D4").Value = 31783

# Row 5 - End Date: 32834 (11/22/1989) -> 32133 (12/22/1987)
$ws.Range("B5").Value = 32133
$ws.Range("C5").Value = 32133
$ws.Range("D5").Value = 32133
